$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.423.65'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.553.84'
$ws.Range('E3').Value = '  -1.92%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.07'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '1.776.12'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').Value = '1.551.03'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range('D14').Value = '28.437.14'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').Value = '0.0₃0673'
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  -3.15%  '
$ws.Range('E30').Value = '  -3.23%  '
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('D33').Value = '1.384.15'
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('E34').Value = '  -3.08%  '
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  -2.64%  '
$ws.Range('E38').Value = '  -2.75%  '
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.509'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.773'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').Value = '1.688.45'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.67%  '
$ws.Range('E51').Value = '  -1.71%  '
